# Add "percentage" and "multiplier" columns to the product import template.
#
# The sheet had a block of header/value columns starting at N (attribute_ids,
# attribute_texts, categories, SEO). Two new columns - "percentage" and
# "multiplier" - are inserted in front of that block (new N:O), pushing the
# existing N:Q block to P:S. Only row 3 gets values for the new columns
# (0.8 and 2); row 2 stays blank for them, matching the rest of that block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new, blank columns at N:O - shifts old N:Q (and their column
# width formatting) two columns to the right, to P:S.
$ws.Columns("N:O").Insert()

# New header row cells for the inserted columns.
$ws.Range("N1").Value = "percentage"
$ws.Range("O1").Value = "multiplier"

# New data values (row 3 only).
$ws.Range("N3").Value = 0.8
$ws.Range("O3").Value = 2

# Give the new columns the same custom width as their neighbor column M.
$ws.Range("N:O").ColumnWidth = $ws.Range("M1").ColumnWidth

# Restore/update the view state: scrolled so column G is left-most visible,
# with O9 as the active selected cell.
$ws.Range("O9").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1 | Out-Null
